$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Config")
$ws1.Hyperlinks.Add($ws1.Range("C6"), "http://kliappsitsvr/OnlineTermInsuranceNRI/") | Out-Null
$ws1.Range("C4").Copy($ws1.Range("C6")) | Out-Null
$ws1.Range("C6").Value = "http://kliappsitsvr/OnlineTermInsuranceNRI/"
